$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) column cells are treated as text so numeric-looking strings
# (e.g. "326.34", "0.4590") are not auto-converted to numbers by Excel.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.896.15'
$ws.Range('E2').Value = '  +1.51%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.893.87'
$ws.Range('E3').Value = '  +1.60%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '326.34'
$ws.Range('E5').Value = '  +0.62%  '

$ws.Range('E6').Value = '  -0.28%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4590'
$ws.Range('E7').Value = '  +0.93%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3903'
$ws.Range('E8').Value = '  +2.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07840'
$ws.Range('E9').Value = '  +0.34%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9888'
$ws.Range('E10').Value = '  +0.25%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.91'
$ws.Range('E11').Value = '  +1.87%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.852.84'
$ws.Range('E12').Value = '  -0.36%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.040'
$ws.Range('E13').Value = '  +2.13%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.700'
$ws.Range('E14').Value = '  +1.59%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06964'
$ws.Range('E15').Value = '  +0.85%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.06'
$ws.Range('E16').Value = '  +1.65%  '

$ws.Range('E17').Value = '  -0.31%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009956'
$ws.Range('E18').Value = '  +0.33%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.92'
$ws.Range('E19').Value = '  +1.92%  '

$ws.Range('E20').Value = '  -0.16%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '28.888.83'
$ws.Range('E21').Value = '  +1.47%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.294'
$ws.Range('E22').Value = '  +1.10%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.98'
$ws.Range('E23').Value = '  +0.98%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.134.04'
$ws.Range('E24').Value = '  +2.49%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.064'
$ws.Range('E25').Value = '  -1.55%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '156.02'
$ws.Range('E26').Value = '  +1.54%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.32'
$ws.Range('E27').Value = '  +1.28%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.929'
$ws.Range('E28').Value = '  +4.92%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.931'
$ws.Range('E29').Value = '  +0.78%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '117.66'
$ws.Range('E30').Value = '  +0.36%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09349'
$ws.Range('E31').Value = '  +1.03%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9130'
$ws.Range('E32').Value = '  +1.10%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.292'
$ws.Range('E33').Value = '  +0.88%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.335'
$ws.Range('E34').Value = '  +1.30%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.282'
$ws.Range('E35').Value = '  -0.34%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.190'
$ws.Range('E36').Value = '  +4.48%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05764'
$ws.Range('E37').Value = '  +1.36%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02073'
$ws.Range('E38').Value = '  +1.08%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.002'
$ws.Range('E39').Value = '  -0.31%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.776'
$ws.Range('E40').Value = '  +1.73%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5679'
$ws.Range('E41').Value = '  +2.43%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1772'
$ws.Range('E42').Value = '  +0.35%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '9.782'
$ws.Range('E43').Value = '  +1.83%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.252'
$ws.Range('E44').Value = '  +7.14%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '11.88'
$ws.Range('E45').Value = '  +2.84%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5350'
$ws.Range('E46').Value = '  +2.37%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.07046'
$ws.Range('E47').Value = '  -0.52%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.841'
$ws.Range('E48').Value = '  +1.93%  '

$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.529'
$ws.Range('E49').Value = '  +4.25%  '

$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '112.32'
$ws.Range('E50').Value = '  +0.49%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.068'
$ws.Range('E51').Value = '  -5.57%  '
